$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G3 and H3 go from 0 to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: H4 goes from 0 to 1
$ws.Range("H4").Value = 1

# Row 5: H5 goes from 0 to 1
$ws.Range("H5").Value = 1

# Row 6: D6 and E6 go from 0 to 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-18: H column goes from 0 to 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
